# Manage Interviewers - Bug Fixes
#
# LIVE_INTERVIEW_HISTORY_DATA.xlsx / "AMS" sheet:
#  - Row 12 (2021-06-16, live_145_hf2): the captured run time had a
#    sub-second rounding issue - correct it.
#  - Rows 13 & 14 had been left as empty placeholder rows; fill them in
#    with the two sprint runs that were actually executed on
#    2021-06-17 (145_data_hstry and 145_hstry_data).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AMS")

# ---------------------------------------------------------------------
# Row 12 - fix the recorded run time (everything else is already right)
# ---------------------------------------------------------------------
$ws.Range("B12").Value2 = 44363.79865090278

# ---------------------------------------------------------------------
# Row 13 - new row: 2021-06-17 / 145_data_hstry
# ---------------------------------------------------------------------
# Column A holds a plain text date label ("2021-06-17"), not a real
# date value, so stage it through a helper cell formatted as Text and
# paste just the value across - this avoids Excel's automatic
# text->date conversion while typing it directly into A13 would trigger.
$helper = $ws.Range("ZZ1")
$helper.NumberFormat = "@"
$helper.Value = "2021-06-17"
$helper.Copy()
$ws.Range("A13").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$helper.Clear()

$ws.Range("B13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B13").Value2 = 44364.553549375
$ws.Range("C13").Value = "145_data_hstry"
$ws.Range("D13").Value = 105
$ws.Range("E13").Value = 105
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 2.87

# ---------------------------------------------------------------------
# Row 14 - new row: 2021-06-17 / 145_hstry_data
# ---------------------------------------------------------------------
$helper = $ws.Range("ZZ1")
$helper.NumberFormat = "@"
$helper.Value = "2021-06-17"
$helper.Copy()
$ws.Range("A14").PasteSpecial(-4163)
$helper.Clear()

$ws.Range("B14").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B14").Value2 = 44364.61715227317
$ws.Range("C14").Value = "145_hstry_data"
$ws.Range("D14").Value = 105
$ws.Range("E14").Value = 105
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 2.89
